$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ranking (A), business_status (B), name (C), rating (E), and user_ratings_total (F)
# for rows 2-54 to reflect the refreshed scrape data, and merge row 55 into row 54
# before removing the now-redundant last row.

$ws.Cells.Item(2, 1).Value = 17
$ws.Cells.Item(2, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(2, 3).Value = 'Arizona Sport Shirts'
$ws.Cells.Item(2, 5).Value = 5
$ws.Cells.Item(2, 6).Value = 38

$ws.Cells.Item(3, 1).Value = 15
$ws.Cells.Item(3, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(3, 3).Value = 'Avon Sports Apparel Corporation'
$ws.Cells.Item(3, 5).Value = 4.7
$ws.Cells.Item(3, 6).Value = 54

$ws.Cells.Item(4, 1).Value = 11
$ws.Cells.Item(4, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(4, 3).Value = 'Champs Sports'
$ws.Cells.Item(4, 4).Value = 2
$ws.Cells.Item(4, 5).Value = 4.2
$ws.Cells.Item(4, 6).Value = 55

$ws.Cells.Item(5, 1).Value = 8
$ws.Cells.Item(5, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(5, 3).Value = 'DICK''S Sporting Goods'
$ws.Cells.Item(5, 5).Value = 4
$ws.Cells.Item(5, 6).Value = 384

$ws.Cells.Item(6, 1).Value = 46
$ws.Cells.Item(6, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(6, 3).Value = 'Elevation Golf Indianapolis | Hamilton County Sports'
$ws.Cells.Item(6, 5).Value = 5
$ws.Cells.Item(6, 6).Value = 5

$ws.Cells.Item(7, 1).Value = 48
$ws.Cells.Item(7, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(7, 3).Value = 'Finish Line'
$ws.Cells.Item(7, 4).Value = 2
$ws.Cells.Item(7, 5).Value = 4.3
$ws.Cells.Item(7, 6).Value = 123

$ws.Cells.Item(8, 1).Value = 41
$ws.Cells.Item(8, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(8, 3).Value = 'Genesis Sports Performance'
$ws.Cells.Item(8, 5).Value = 4.9
$ws.Cells.Item(8, 6).Value = 9

$ws.Cells.Item(9, 1).Value = 9
$ws.Cells.Item(9, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(9, 3).Value = 'Hibbett Sports'
$ws.Cells.Item(9, 5).Value = 3.9
$ws.Cells.Item(9, 6).Value = 97

$ws.Cells.Item(10, 1).Value = 51
$ws.Cells.Item(10, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(10, 3).Value = 'Indiana Sports Corp'
$ws.Cells.Item(10, 5).Value = 5
$ws.Cells.Item(10, 6).Value = 2

$ws.Cells.Item(11, 1).Value = 31
$ws.Cells.Item(11, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(11, 3).Value = 'Indy Indoor Sport'
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = 0

$ws.Cells.Item(12, 1).Value = 58
$ws.Cells.Item(12, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(12, 3).Value = 'Indy Sport'
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 6).Value = 0

$ws.Cells.Item(13, 1).Value = 35
$ws.Cells.Item(13, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(13, 3).Value = 'Indy Sport Group'
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0

$ws.Cells.Item(14, 1).Value = 20
$ws.Cells.Item(14, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(14, 3).Value = 'Indy Sports Performance'
$ws.Cells.Item(14, 5).Value = 5
$ws.Cells.Item(14, 6).Value = 2

$ws.Cells.Item(15, 1).Value = 25
$ws.Cells.Item(15, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(15, 3).Value = 'KS&E Sports'
$ws.Cells.Item(15, 5).Value = 2.8
$ws.Cells.Item(15, 6).Value = 34

$ws.Cells.Item(16, 1).Value = 55
$ws.Cells.Item(16, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(16, 3).Value = 'Larson and Sons Target Sports'
$ws.Cells.Item(16, 5).Value = 4.4
$ws.Cells.Item(16, 6).Value = 52

$ws.Cells.Item(17, 1).Value = 21
$ws.Cells.Item(17, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(17, 3).Value = 'Looking Good Sports Llc'
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = 0

$ws.Cells.Item(18, 1).Value = 33
$ws.Cells.Item(18, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(18, 3).Value = 'Marksman Shooting Sports'
$ws.Cells.Item(18, 5).Value = 4.9
$ws.Cells.Item(18, 6).Value = 51

$ws.Cells.Item(19, 1).Value = 42
$ws.Cells.Item(19, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(19, 3).Value = 'Midwest Sports Complex'
$ws.Cells.Item(19, 5).Value = 3.5
$ws.Cells.Item(19, 6).Value = 284

$ws.Cells.Item(20, 1).Value = 3
$ws.Cells.Item(20, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(20, 3).Value = 'N & D Sports'
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = 0

$ws.Cells.Item(21, 1).Value = 24
$ws.Cells.Item(21, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(21, 3).Value = 'National Institute for Fitness and Sport (NIFS)'
$ws.Cells.Item(21, 5).Value = 4.3
$ws.Cells.Item(21, 6).Value = 52

$ws.Cells.Item(22, 1).Value = 40
$ws.Cells.Item(22, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(22, 3).Value = 'Oakley Store'
$ws.Cells.Item(22, 5).Value = 4.6
$ws.Cells.Item(22, 6).Value = 78

$ws.Cells.Item(23, 1).Value = 59
$ws.Cells.Item(23, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(23, 3).Value = 'Off The Wall Sports'
$ws.Cells.Item(23, 5).Value = 4.6
$ws.Cells.Item(23, 6).Value = 234

$ws.Cells.Item(24, 1).Value = 19
$ws.Cells.Item(24, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(24, 3).Value = 'Performance Sports'
$ws.Cells.Item(24, 5).Value = 3.3
$ws.Cells.Item(24, 6).Value = 3

$ws.Cells.Item(25, 1).Value = 2
$ws.Cells.Item(25, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(25, 3).Value = 'Play It Again Sports'
$ws.Cells.Item(25, 5).Value = 4.2
$ws.Cells.Item(25, 6).Value = 103

$ws.Cells.Item(26, 1).Value = 18
$ws.Cells.Item(26, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(26, 3).Value = 'Ruben Sports'
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(26, 6).Value = 0

$ws.Cells.Item(27, 1).Value = 29
$ws.Cells.Item(27, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(27, 3).Value = 'Sport Clips Haircuts of Avon'
$ws.Cells.Item(27, 5).Value = 4
$ws.Cells.Item(27, 6).Value = 128

$ws.Cells.Item(28, 1).Value = 50
$ws.Cells.Item(28, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(28, 3).Value = 'Sport Clips Haircuts of Carmel'
$ws.Cells.Item(28, 5).Value = 4
$ws.Cells.Item(28, 6).Value = 97

$ws.Cells.Item(29, 1).Value = 30
$ws.Cells.Item(29, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(29, 3).Value = 'Sport Clips Haircuts of Carmel - 146th Street'
$ws.Cells.Item(29, 5).Value = 3.9
$ws.Cells.Item(29, 6).Value = 57

$ws.Cells.Item(30, 1).Value = 37
$ws.Cells.Item(30, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(30, 3).Value = 'Sport Clips Haircuts of Carmel - Carmel Point'
$ws.Cells.Item(30, 5).Value = 4.6
$ws.Cells.Item(30, 6).Value = 158

$ws.Cells.Item(31, 1).Value = 57
$ws.Cells.Item(31, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(31, 3).Value = 'Sport Clips Haircuts of Carmel - E 126th'
$ws.Cells.Item(31, 5).Value = 4.6
$ws.Cells.Item(31, 6).Value = 54

$ws.Cells.Item(32, 1).Value = 34
$ws.Cells.Item(32, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(32, 3).Value = 'Sport Clips Haircuts of Castleton Crossing'
$ws.Cells.Item(32, 5).Value = 4.4
$ws.Cells.Item(32, 6).Value = 126

$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(33, 3).Value = 'Sport Clips Haircuts of Fishers'
$ws.Cells.Item(33, 5).Value = 4
$ws.Cells.Item(33, 6).Value = 108

$ws.Cells.Item(34, 1).Value = 49
$ws.Cells.Item(34, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(34, 3).Value = 'Sport Clips Haircuts of Indianapolis - 96th Street'
$ws.Cells.Item(34, 5).Value = 4.3
$ws.Cells.Item(34, 6).Value = 112

$ws.Cells.Item(35, 1).Value = 38
$ws.Cells.Item(35, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(35, 3).Value = 'Sport Clips Haircuts of Medford Place'
$ws.Cells.Item(35, 5).Value = 4.3
$ws.Cells.Item(35, 6).Value = 101

$ws.Cells.Item(36, 1).Value = 45
$ws.Cells.Item(36, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(36, 3).Value = 'Sport Clips Haircuts of Noblesville at Stoney Creek Commons'
$ws.Cells.Item(36, 5).Value = 4.6
$ws.Cells.Item(36, 6).Value = 261

$ws.Cells.Item(37, 1).Value = 27
$ws.Cells.Item(37, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(37, 3).Value = 'Sport Clips Haircuts of North Keystone'
$ws.Cells.Item(37, 5).Value = 3.9
$ws.Cells.Item(37, 6).Value = 100

$ws.Cells.Item(38, 1).Value = 56
$ws.Cells.Item(38, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(38, 3).Value = 'Sport Clips Haircuts of Northfield Commons'
$ws.Cells.Item(38, 5).Value = 4.5
$ws.Cells.Item(38, 6).Value = 148

$ws.Cells.Item(39, 1).Value = 36
$ws.Cells.Item(39, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(39, 3).Value = 'Sport Clips Haircuts of Shadeland Place'
$ws.Cells.Item(39, 5).Value = 4.3
$ws.Cells.Item(39, 6).Value = 91

$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(40, 3).Value = 'Sport Clips Haircuts of Westfield'
$ws.Cells.Item(40, 5).Value = 3.7
$ws.Cells.Item(40, 6).Value = 65

$ws.Cells.Item(41, 1).Value = 52
$ws.Cells.Item(41, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(41, 3).Value = 'Sport Clips Haircuts of Zionsville'
$ws.Cells.Item(41, 5).Value = 4.2
$ws.Cells.Item(41, 6).Value = 96

$ws.Cells.Item(42, 1).Value = 26
$ws.Cells.Item(42, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(42, 3).Value = 'Sport Graphics Inc'
$ws.Cells.Item(42, 5).Value = 4.3
$ws.Cells.Item(42, 6).Value = 10

$ws.Cells.Item(43, 1).Value = 1
$ws.Cells.Item(43, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(43, 3).Value = 'Sport Passes'
$ws.Cells.Item(43, 5).Value = 0
$ws.Cells.Item(43, 6).Value = 0

$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(44, 3).Value = 'Sport''n Image'
$ws.Cells.Item(44, 5).Value = 5
$ws.Cells.Item(44, 6).Value = 3

$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(45, 3).Value = 'Sport.ly'
$ws.Cells.Item(45, 5).Value = 0
$ws.Cells.Item(45, 6).Value = 0

$ws.Cells.Item(46, 1).Value = 16
$ws.Cells.Item(46, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(46, 3).Value = 'Sports Corporation Inc'
$ws.Cells.Item(46, 5).Value = 0
$ws.Cells.Item(46, 6).Value = 0

$ws.Cells.Item(47, 1).Value = 54
$ws.Cells.Item(47, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(47, 3).Value = 'Sports Select'
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(47, 6).Value = 0

$ws.Cells.Item(48, 1).Value = 0
$ws.Cells.Item(48, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(48, 3).Value = 'Sports Spot'
$ws.Cells.Item(48, 5).Value = 3.3
$ws.Cells.Item(48, 6).Value = 13

$ws.Cells.Item(49, 1).Value = 28
$ws.Cells.Item(49, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(49, 3).Value = 'Sports Travel & Tickets'
$ws.Cells.Item(49, 5).Value = 0
$ws.Cells.Item(49, 6).Value = 0

$ws.Cells.Item(50, 1).Value = 53
$ws.Cells.Item(50, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(50, 3).Value = 'St. Vincent Sports Performance'
$ws.Cells.Item(50, 5).Value = 5
$ws.Cells.Item(50, 6).Value = 1

$ws.Cells.Item(51, 1).Value = 13
$ws.Cells.Item(51, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(51, 3).Value = 'Sun Valley Sports'
$ws.Cells.Item(51, 5).Value = 4.5
$ws.Cells.Item(51, 6).Value = 112

$ws.Cells.Item(52, 1).Value = 14
$ws.Cells.Item(52, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(52, 3).Value = 'Team Sports'
$ws.Cells.Item(52, 5).Value = 5
$ws.Cells.Item(52, 6).Value = 9

$ws.Cells.Item(53, 1).Value = 22
$ws.Cells.Item(53, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(53, 3).Value = 'The North Face The Fashion Mall at Keystone'
$ws.Cells.Item(53, 5).Value = 4.2
$ws.Cells.Item(53, 6).Value = 124

$ws.Cells.Item(54, 1).Value = 10
$ws.Cells.Item(54, 2).Value = 'OPERATIONAL'
$ws.Cells.Item(54, 3).Value = 'Webster''s Sporting Goods'
$ws.Cells.Item(54, 5).Value = 4.5
$ws.Cells.Item(54, 6).Value = 28

# Row 55 (Webster's Sporting Goods duplicate source row) is no longer needed
# now that its data has been merged into row 54; remove it so the table ends at row 54.
$ws.Rows.Item(55).Delete()
